$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.539.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.294.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.644.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.296.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.498.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +35.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.04%  "
$ws.Range("E25").Value = "  -2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0879"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.88%  "
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "80.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.697.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.90%  "
